$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update Random Forest row (row 2): Precision, Recall, F1-Score
$ws.Range("B2").Value = 0.2171428571428571
$ws.Range("C2").Value = 0.06280991735537191
$ws.Range("D2").Value = 0.09743589743589744

# Update LGBM row (row 3): Precision, Recall, F1-Score
$ws.Range("B3").Value = 0.1612903225806452
$ws.Range("C3").Value = 0.03305785123966942
$ws.Range("D3").Value = 0.05486968449931413
